$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark (currently sits at the end
#        of the "Writer: Hoang Trung Hieu" paragraph). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Locate the "Reviewer:" paragraph and replace its trailing name
#        run ("Hoang Trung Hieu") with "Mr. Kieu Trong Khanh.", keeping
#        the run split from the preceding " " run intact (InsertXML
#        swaps in a standalone run instead of a plain text assignment,
#        which would otherwise silently merge it into the neighbour). ---
$reviewerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Reviewer:")) {
        $reviewerPara = $p
    }
}

$probeRange = $reviewerPara.Range.Duplicate
$probeRange.MoveEnd(1, -1)
$probeRange.MoveStart(1, $probeRange.Text.Length - "Hoang Trung Hieu".Length)
$nameStart = $probeRange.Start
$nameEnd = $probeRange.End

# Re-derive a fresh Range from plain numeric bounds - InsertXML only
# replaces the targeted span when it is handed a range built straight
# from Document.Range(start, end); a range that was produced via
# Duplicate/Move degrades to a plain insert instead of a replace.
$nameRange = $d.Range($nameStart, $nameEnd)
$runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Mr. Kieu Trong Khanh.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$nameRange.InsertXML($runXml)

# --- 3. Re-add the "_GoBack" bookmark right after the new text, at the
#        end of the "Reviewer:" paragraph (before its paragraph mark).
#        A bookmark can't be inserted collapsed exactly at
#        paragraph.End-1 reliably, so a throw-away marker character is
#        appended, the bookmark is anchored just before it, and the
#        marker is removed again - leaving the bookmark collapsed right
#        after the last run, matching the original placement style. ---
$reviewerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Reviewer:")) {
        $reviewerPara = $p
    }
}

$endRange = $d.Range($reviewerPara.Range.End - 1, $reviewerPara.Range.End - 1)
$endRange.InsertAfter("X")

$reviewerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Reviewer:")) {
        $reviewerPara = $p
    }
}

$bookmarkPos = $reviewerPara.Range.End - 2
$target = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $target)

$markerRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$markerRange.Delete()
